# Modificado modulo de inicio
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 9: move the "X" marker from column E ("Por hacer") to column G ("Terminado")
$ws.Range("E9").Value = $null
$ws.Range("G9").Value = "X"

# Row 10: mark column E ("Por hacer") with "X"
$ws.Range("E10").Value = "X"

# Row 11: mark column E ("Por hacer") with "X"
$ws.Range("E11").Value = "X"

# Update the active selection to E12, matching the recorded cursor position
$ws.Range("E12").Select()
